$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.271.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.632.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.521'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.46%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0850'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.617.75'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.545'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.193.94'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0₃0735'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '217.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.316.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.52%  '
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.849'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.541'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.799'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.770.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('E48').Value = '  +21.42%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0963'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.78%  '
